# BOM.xlsx update:
#  - Add more components to the BOM (IEC connector, IEC->UK power cord,
#    cartridge fuse, second colour of 18awg wire)
#  - Record a couple of "Brought Myself" prices and a batch-1 price that
#    were missing
#
# All numeric-looking values in this sheet are stored as TEXT (the sheet
# uses t="str" cells for prices), so every such value below is entered
# with a leading apostrophe to force Excel to keep it as text instead of
# silently converting it to a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- first do all the row insertions, THEN fill in values -----------------
# (inserting rows after a cell has been given a "quote-prefix" text style
#  can leak that style into the freshly inserted row, so shift everything
#  into place first and only write values afterwards)

# 3 new rows for new Cables-category components: pushes the existing
# "Cables" group header + everything after it down by 3, from row 17 to
# row 20.
$ws.Rows("17:19").Insert()

# 1 more new row for the second colour of 18awg wire, right after the
# (now shifted) "18 awg" row which lives at row 25.
$ws.Rows("26:26").Insert()

# --- fill in a couple of prices that already had a row -------------------

# Okdo OV5647 Adjustable Focus now has a recorded batch-1 purchase price
$ws.Range("J10").Value = "'9.83"

# Ribbon Extension Cable + LED Ring : "Brought Myself" prices
$ws.Range("M12").Value = "'2.20"
$ws.Range("M14").Value = "'6.00"

# --- new Cables-category rows (17-19) --------------------------------------

$ws.Range("D17").Value = "C14 IEC20 Connector"
$ws.Range("E17").Value = "https://uk.rs-online.com/web/p/iec-connectors/8117216"
$ws.Range("F17").Value = "'3.71"
$ws.Range("G17").Value = "2pc"
$ws.Range("K17").Value = "'7.42"

$ws.Range("D18").Value = "IEC C13 Socket to Type G UK Plug"
$ws.Range("E18").Value = "https://uk.rs-online.com/web/p/power-cords/1469109?gb=s"
$ws.Range("F18").Value = "'4.46"
$ws.Range("G18").Value = "1.8m"
$ws.Range("K18").Value = "'4.46"

$ws.Range("D19").Value = "6A T Glass Cartridge Fuse, 5 x 20mm"
$ws.Range("E19").Value = "https://uk.rs-online.com/web/p/cartridge-fuses/9113427?gb=s"
$ws.Range("F19").Value = "'0.333"
$ws.Range("G19").Value = "10pc"
$ws.Range("K19").Value = "'3.33"

# --- rename the existing wire row to be colour-specific & fix its price ---
$ws.Range("D25").Value = "18 awg red"
$ws.Range("J25").Value = "'10.27"

# --- new row for the black 18awg wire --------------------------------------
$ws.Range("D26").Value = "18 awg black"
$ws.Range("E26").Value = "https://uk.rs-online.com/web/p/hook-up-wire/8114416?gb=s"
$ws.Range("F26").Value = "'10.27"
$ws.Range("G26").Value = "25m"
$ws.Range("K26").Value = "'10.27"
